# add ACME dev requirements
$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Mark the "Project management plan" inline picture as NoProofing
#    (adds <w:rPr><w:noProof/></w:rPr> to the run that holds the <w:drawing>)
#    Identify it by its known size (5731510 x 1027430 EMU = 451.3 x 80.9 pt)
#    rather than by position, so only that single picture's run is touched.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $shp = $d.InlineShapes($i)
    if ([Math]::Abs($shp.Width - 451.3) -lt 0.5 -and [Math]::Abs($shp.Height - 80.9) -lt 0.5) {
        $shp.Range.NoProofing = $true
    }
}

# ---------------------------------------------------------------------------
# 2) Replace the empty "Acme entertainment development requirements" body
#    paragraph (ind left=432, holds the _GoBack bookmark) with the five new
#    ACME dev-requirements bullet points + a trailing blank ind-432 paragraph.
# ---------------------------------------------------------------------------
$bodyPara = $null
$pCount = $d.Paragraphs.Count
for ($i = 1; $i -le $pCount; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -match "^[\r\a]?$") {
        $prev = $null
        if ($i -gt 1) { $prev = $d.Paragraphs($i - 1) }
        if ($prev -ne $null -and $prev.Range.Text -match "Acme entertainment development requirements") {
            $bodyPara = $para
            break
        }
    }
}

if ($bodyPara -ne $null) {
    $xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
        '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Deploy a movie database application</w:t></w:r></w:p>' + `
        '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>RAD agile methodology</w:t></w:r></w:p>' + `
        '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Movie website must be able to search</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>' + `
        '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Movie website must be able to display top 10</w:t></w:r></w:p>' + `
        '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Movie website must be responsive to different screen sizes</w:t></w:r></w:p>' + `
        '<w:p><w:pPr><w:ind w:left="432"/></w:pPr></w:p>' + `
        '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    [void]$bodyPara.Range.InsertXML($xml)
}

Write-Output "done"
